# Updates the single-column results table:
#  - rows 1-3 get their memory-usage values rewritten to "0M"
#  - 10 new rows are inserted right after (old) row 3, holding the
#    per-iteration timing values that used to be packed (with tabs)
#    into the last three rows of the table
#  - the last three rows, which held 10 tab-separated values each,
#    are collapsed down to just their leading "iteration count" value

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) first three single-value rows -> "0M"
$t.Cell(1, 1).Range.Paragraphs.Item(1).Range.Text = "0M"
$t.Cell(2, 1).Range.Paragraphs.Item(1).Range.Text = "0M"
$t.Cell(3, 1).Range.Paragraphs.Item(1).Range.Text = "0M"

# 2) insert 10 new rows right after row 3, each with a single value
$newValues = @("104", "0.00002", "0.00009", "0.00006", "0.00001", `
                "0.00008", "0.00009", "0.00009", "0.00452", "100.0")

$insertPos = 4
foreach ($val in $newValues) {
    $refRow = $t.Rows.Item($insertPos)
    $t.Rows.Add($refRow) | Out-Null
    $t.Cell($insertPos, 1).Range.Paragraphs.Item(1).Range.Text = $val
    $insertPos = $insertPos + 1
}

# 3) the final three rows (now shifted down by the 10 inserted rows)
#    had their tab-separated values collapsed to just the first value
$lastCount = $t.Rows.Count
$t.Cell($lastCount - 2, 1).Range.Paragraphs.Item(1).Range.Text = "100"
$t.Cell($lastCount - 1, 1).Range.Paragraphs.Item(1).Range.Text = "0"
$t.Cell($lastCount, 1).Range.Paragraphs.Item(1).Range.Text = "206"
